# Insert a new data row at row 60 (pushing the existing rows 60:131 down to
# 61:132), matching the weekly refresh that prepends the newest price
# observation to the "Cereza" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 60:131 down by one row.
$ws.Rows(60).Insert()

# Populate the newly inserted row 60 with the new weekly record. Columns not
# specific to this record (mercado/region/categoria metadata) mirror the
# surrounding rows, same as every other row in this table.
$ws.Range("A60").Value = 4
$ws.Range("B60").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C60").Value = "Los Lagos"
$ws.Range("D60").Value = 45280
$ws.Range("E60").Value = 10
$ws.Range("F60").Value = "Fruta"
$ws.Range("G60").Value = 100103
$ws.Range("H60").Value = "Frutos de hueso (carozo)"
$ws.Range("I60").Value = 100103001
$ws.Range("J60").Value = "Cereza"
$ws.Range("K60").Value = "Lapins"
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 350
$ws.Range("N60").Value = 11000
$ws.Range("O60").Value = 11000
$ws.Range("P60").Value = 11000
$ws.Range("Q60").Value = "$/bandeja 10 kilos"
$ws.Range("R60").Value = "Provincia de Curicó"
$ws.Range("S60").Value = 1100
$ws.Range("T60").Value = 10
